# Update res_bus/vm_pu.xlsx values: slack bus voltage setpoint changed from 1.05 to 1.02 (380 kV case),
# with recalculated bus voltage magnitudes (p.u.) for data rows 2-25 (index 0-23) across columns
# B,C,D,E,F,I,J,K,L,M,N. Columns G (constant 1) and H (blank) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column letters (in order) whose values change for every data row
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New values per data row (rows 2..25), one array per row aligned with $cols
$newValues = @(
    @(1.02,1.070034384315117,1.048162697484104,1.073888080771541,1.083075494277823,1.052767953902756,1.074966085611274,1.050923311132038,1.076578749374504,1.085742052965537,1.076492660753175),
    @(1.02,1.07157280831684,1.048880437935124,1.075267929392525,1.084546861236421,1.053212776556474,1.076159124924949,1.051453285663418,1.077774235732329,1.087030596296606,1.077687394319519),
    @(1.02,1.072566765406819,1.049344120932974,1.076159517611225,1.085497845376303,1.053498617270882,1.076929168365494,1.051794801590235,1.078545992587246,1.087862759373424,1.078458531310075),
    @(1.02,1.072984272344975,1.049538876905232,1.076534044555907,1.085897385178928,1.05361831074529,1.077252437928068,1.051938038410369,1.07887001346725,1.088212221296908,1.078782259952507),
    @(1.02,1.073054353149113,1.049571566992251,1.076596912061144,1.085964454973303,1.053638380078841,1.077306689702969,1.051962068823876,1.078924393205727,1.088270875419499,1.078836588771151),
    @(1.02,1.072572345530699,1.049346723962912,1.076164523216712,1.085503185036182,1.053500218480488,1.076933489696309,1.051796716847534,1.078550323835555,1.087867430384064,1.078462858777675),
    @(1.02,1.070554618552743,1.04840541601528,1.074354673071131,1.083572979090052,1.052918697071807,1.075369682653479,1.051102712177098,1.076983146409438,1.08617785941975,1.076896830949474),
    @(1.02,1.066987233941491,1.046740950377116,1.071155502837127,1.080163070516155,1.051878643819314,1.072598981895755,1.049868883349711,1.074207508128173,1.083188007487887,1.07412219547887),
    @(1.02,1.064600481217922,1.045627333005545,1.069015608383513,1.077883566394281,1.051174817347994,1.070741333950553,1.049038883255621,1.072347250020687,1.081185904384812,1.072261909460511),
    @(1.02,1.063564863713454,1.045144158687211,1.068087232899305,1.076894943012253,1.050867539549278,1.069934372918841,1.048677691083048,1.071539322242764,1.080316780984755,1.071453802451574),
    @(1.02,1.063179858349864,1.04496453830931,1.067742116573449,1.076527478179659,1.050753021956911,1.069634235538956,1.048543255939511,1.071238850490606,1.079993611803902,1.071153238842183),
    @(1.02,1.063262458438692,1.045003074217847,1.067816157828278,1.076606311902935,1.050777603654084,1.069698634070974,1.048572105093961,1.071303319625362,1.08006294811038,1.071217728827503),
    @(1.02,1.063533045857022,1.045129314240536,1.06805871115892,1.076864573322694,1.050858081277174,1.069909571593323,1.0486665842005,1.071514492772705,1.080290074659765,1.071428965905329),
    @(1.02,1.063699719554725,1.045207075220629,1.068208119418286,1.077023663849785,1.050907615632147,1.070039484449488,1.048724759744294,1.071644554085733,1.080429969799516,1.071559063252651),
    @(1.02,1.064669165870817,1.04565937906526,1.069077183286899,1.077949144015974,1.051195157101845,1.070794834154664,1.049062816336106,1.072400817842366,1.081243538253616,1.072315485641048),
    @(1.02,1.065276694731033,1.045942835977822,1.069621840120685,1.078529243636965,1.051374848552201,1.071267947175383,1.049274387629701,1.072874548101954,1.081753274025777,1.072789270536523),
    @(1.02,1.065630850428795,1.046108078089541,1.069939357345002,1.078867453750961,1.051479416946382,1.071543656872336,1.049397620395986,1.073150633585541,1.082050382570138,1.07306537177287),
    @(1.02,1.065751573787819,1.046164405581133,1.070047593593522,1.078982749095725,1.051515031006208,1.071637624721792,1.049439610279528,1.07324473223715,1.082151653225167,1.07315947306745),
    @(1.02,1.065211533971393,1.045912433421098,1.069563421433496,1.078467020238214,1.051355594473224,1.071217212447583,1.049251705963175,1.072823745518809,1.081698606154478,1.072738463759591),
    @(1.02,1.063453373819801,1.045092143800655,1.067987292940063,1.076788528631465,1.050834393162084,1.069847466745301,1.048638769994039,1.071452317855293,1.080223200945144,1.071366772861299),
    @(1.02,1.062346029327974,1.044575538265452,1.066994715678862,1.075731766403614,1.050504487780306,1.068983957694216,1.048251816257428,1.070587894393779,1.079293595459855,1.070502037528316),
    @(1.02,1.062933238799713,1.044849482563921,1.067521054219969,1.0762921141784,1.050679586886262,1.069441940265885,1.048457097989308,1.071046348092753,1.079786585220983,1.070960670487767),
    @(1.02,1.065240977949385,1.045926171317043,1.069589818855018,1.078495136784225,1.051364295312437,1.071240138057205,1.049261955367289,1.072846701740646,1.081723308866761,1.072761421926207),
    @(1.02,1.067910948978508,1.047171947565168,1.07198378930629,1.081045680788667,1.05214935503327,1.073317099618352,1.05018916174903,1.074926782801829,1.083962490330762,1.074841333011002)
)

for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $newValues[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowData[$j]
    }
}
